$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1=14, Q1=15 with the same style as O1 ---
$ws.Cells.Item(1, 15).Copy($ws.Cells.Item(1, 16))   # O1 -> P1 (copies style s="1")
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 15).Copy($ws.Cells.Item(1, 17))   # O1 -> Q1 (copies style s="1")
$ws.Cells.Item(1, 17).Value = 15

# --- Data rows 2-25 ---
# Columns (1-based): I=9, J=10, K=11, L=12, M=13, N=14, O=15, P=16, Q=17
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2    # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1    # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2    # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1    # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2    # P: new
    $ws.Cells.Item($r, 17).Value = 2    # Q: new
}
